$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column E (reviews_count) entirely, shifting F:K left to E:J
$ws.Columns("E").Delete()
